$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Schedule sheet updates (cost/unit-cost recompute for run 170)
$wsSchedule.Range("E2").Value = 340.98722775
$wsSchedule.Range("F2").Value = 7.517355109126984
$wsSchedule.Range("E3").Value = 747.545175
$wsSchedule.Range("F3").Value = 28.2518962585034
$wsSchedule.Range("E4").Value = 13.75552425
$wsSchedule.Range("F4").Value = 0.404336397707231

# Detailed sheet updates (Price column recompute + one historical/forecast flag fix)
$wsDetailed.Range("B32").Value = 0.51
$wsDetailed.Range("B33").Value = -5.51
$wsDetailed.Range("B34").Value = -5.50985
$wsDetailed.Range("B35").Value = 0.51
$wsDetailed.Range("B36").Value = 2.98502
$wsDetailed.Range("B37").Value = 9.82779
$wsDetailed.Range("B38").Value = 9.875310000000001
$wsDetailed.Range("B39").Value = 16.41829
$wsDetailed.Range("B40").Value = 22.72033
$wsDetailed.Range("B41").Value = 55.33036
$wsDetailed.Range("B44").Value = 53.45932
$wsDetailed.Range("B45").Value = 53.58333
$wsDetailed.Range("B46").Value = 53.6642
$wsDetailed.Range("B54").Value = 47.49738
$wsDetailed.Range("B55").Value = 48.13534
$wsDetailed.Range("B56").Value = 48.74444
$wsDetailed.Range("B57").Value = 49.47383
$wsDetailed.Range("B58").Value = 50.93627
$wsDetailed.Range("B59").Value = 57.06
$wsDetailed.Range("B60").Value = 56.98
$wsDetailed.Range("B61").Value = 57.66203
$wsDetailed.Range("B62").Value = 57.55805
$wsDetailed.Range("B64").Value = 26.42822
$wsDetailed.Range("B65").Value = 9.697050000000001
$wsDetailed.Range("B68").Value = 0.7
$wsDetailed.Range("B69").Value = 0.0288
$wsDetailed.Range("B70").Value = -0.95029
$wsDetailed.Range("B71").Value = 0.02888
$wsDetailed.Range("B72").Value = 1.77776
$wsDetailed.Range("B73").Value = 0.6
$wsDetailed.Range("B74").Value = 0.59
$wsDetailed.Range("B75").Value = 0.51
$wsDetailed.Range("B76").Value = 0
$wsDetailed.Range("B77").Value = -4.25198
$wsDetailed.Range("B78").Value = -5.02993
$wsDetailed.Range("B79").Value = -6.24082
$wsDetailed.Range("B80").Value = -5.58973
$wsDetailed.Range("B82").Value = -5.24342
$wsDetailed.Range("B83").Value = -6.66332
$wsDetailed.Range("B85").Value = -4.04159
$wsDetailed.Range("B86").Value = 12.2117
$wsDetailed.Range("B87").Value = 25.73014
$wsDetailed.Range("B92").Value = 55.14679
$wsDetailed.Range("B93").Value = 56.57871

$wsDetailed.Range("C34").Value = "historical"
